$d = $word.ActiveDocument

# Locate the "KEY ACHIEVEMENTS AND IMPACT" section so edits are scoped to it only.
# (Several bullet strings in this section duplicate text used verbatim elsewhere
# in the "PROFESSIONAL EXPERIENCE" section, so a document-wide Find must NOT be
# used -- it would also rewrite those unrelated paragraphs.)

$sectionStartPara = $null
$sectionEndPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text.Trim()
    if ($txt -eq "KEY ACHIEVEMENTS AND IMPACT") {
        $sectionStartPara = $i
    } elseif ($sectionStartPara -ne $null -and $txt -eq "TECHNICAL SKILLS") {
        $sectionEndPara = $i
        break
    }
}

$startPos = $d.Paragraphs.Item($sectionStartPara).Range.Start
$endPos = $d.Paragraphs.Item($sectionEndPara).Range.Start
$sectionRange = $d.Range($startPos, $endPos)

# --- Replace the first four bullets of that section ---

$targets = @(
    @{
        Old = [char]0x2022 + " Delivered `$4.9M additional revenue through continuous testing and optimization, increased conversion rates by 23%"
        New = [char]0x2022 + " Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations"
    },
    @{
        Old = [char]0x2022 + " Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations"
        New = [char]0x2022 + " Real-time collaboration at national scale"
    },
    @{
        Old = [char]0x2022 + " Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from " + [char]0x00B1 + "4.2% to " + [char]0x00B1 + "2.1%"
        New = [char]0x2022 + " Revenue generation: Delivered `$4.9M additional revenue through optimization"
    },
    @{
        Old = [char]0x2022 + " Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis"
        New = [char]0x2022 + " 23% conversion rate improvement"
    }
)

foreach ($t in $targets) {
    # Recompute the end of the scoped range each time: replacement text is a
    # different length than the original, which shifts later character offsets.
    $endPos = $d.Paragraphs.Item($sectionEndPara).Range.Start
    $rng = $d.Range($startPos, $endPos)
    $rng.Find.Execute($t.Old, $true, $true, $false, $false, $false, $true, 1, $false, $t.New, 2)
}

# --- Remove the last two bullets of that same section entirely ---
# ("Discovered systematic race coding errors..." and "Developed longitudinal data
# analysis methods..."), which the diff drops.

$deleteOld1 = [char]0x2022 + " Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%"
$deleteOld2 = [char]0x2022 + " Developed longitudinal data analysis methods using geospatial techniques that improved segmentation accuracy by 34% and survey incidence rates by 28%, reducing polling costs while increasing response quality"

$sectionEndPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text.Trim()
    if ($txt -eq "TECHNICAL SKILLS") {
        $sectionEndPara = $i
        break
    }
}

for ($i = $sectionEndPara - 1; $i -ge $sectionStartPara; $i--) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text.Trim()
    if ($text -eq $deleteOld1 -or $text -eq $deleteOld2) {
        $para.Range.Delete()
    }
}
